$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2945437.2
$ws.Range("J17").Value = 3230374.5
$ws.Range("L17").Value = 9691123.5
$ws.Range("N17").Value = -9691459.5
$ws.Range("H32").Value = 364
$ws.Range("J32").Value = 383.33334
$ws.Range("L32").Value = 383.33334
$ws.Range("N32").Value = -1035.33334
$ws.Range("H106").Value = 7753920
$ws.Range("I106").Value = 12821650
$ws.Range("J106").Value = 3273.8823
$ws.Range("K106").Value = 12821650
$ws.Range("L106").Value = 3273.8823
$ws.Range("M106").Value = -12821019
$ws.Range("N106").Value = -4535.8823
$ws.Range("H141").Value = 1050.174
$ws.Range("I141").Value = 762
$ws.Range("J141").Value = 4076
$ws.Range("K141").Value = 2286
$ws.Range("L141").Value = 12228
$ws.Range("M141").Value = 2894
$ws.Range("N141").Value = -22588.5001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4694.465
$ws.Range("I32").Value = 3542.033
$ws.Range("J32").Value = 17803.375
$ws.Range("K32").Value = 3542.033
$ws.Range("L32").Value = 17803.375
$ws.Range("M32").Value = -3255.033
$ws.Range("N32").Value = -18377.375
$ws.Range("H60").Value = 20500
$ws.Range("J60").Value = 20500
$ws.Range("L60").Value = 20500
$ws.Range("N60").Value = -21966
$ws.Range("H61").Value = 360889.03
$ws.Range("I61").Value = 439974.34
$ws.Range("J61").Value = 611.44446
$ws.Range("K61").Value = 439974.34
$ws.Range("L61").Value = 611.44446
$ws.Range("M61").Value = -439762.34
$ws.Range("N61").Value = -1035.44446
$ws.Range("H102").Value = 1650
$ws.Range("I102").Value = 1300
$ws.Range("K102").Value = 1300
$ws.Range("M102").Value = 322
$ws.Range("H132").Value = 11320.82
$ws.Range("I132").Value = 1283.6666
$ws.Range("J132").Value = 64015.875
$ws.Range("K132").Value = 3850.9998
$ws.Range("L132").Value = 192047.625
$ws.Range("M132").Value = -1320.9998
$ws.Range("N132").Value = -197107.625
$ws.Range("H136").Value = 360889.03
$ws.Range("I136").Value = 439974.34
$ws.Range("J136").Value = 611.44446
$ws.Range("K136").Value = 1319923.02
$ws.Range("L136").Value = 1834.33338
$ws.Range("M136").Value = -1317373.02
$ws.Range("N136").Value = -6934.33338

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("N43").Value = 0
$ws.Range("H134").Value = 3764.4324
$ws.Range("J134").Value = 2562.5
$ws.Range("L134").Value = 7687.5
$ws.Range("N134").Value = -12757.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2946.5908
$ws.Range("I31").Value = 1660.6061
$ws.Range("K31").Value = 1660.6061
$ws.Range("M31").Value = -1365.6061
$ws.Range("H34").Value = 2946.5908
$ws.Range("I34").Value = 1660.6061
$ws.Range("K34").Value = 1660.6061
$ws.Range("M34").Value = -1458.6061
$ws.Range("H58").Value = 13987.632
$ws.Range("I58").Value = 838.7059
$ws.Range("K58").Value = 838.7059
$ws.Range("M58").Value = -635.7059
$ws.Range("H132").Value = 1621.7843
$ws.Range("I132").Value = 1233.94
$ws.Range("K132").Value = 3701.82
$ws.Range("M132").Value = -1171.82
$ws.Range("H134").Value = 740.10767
$ws.Range("I134").Value = 623.88135
$ws.Range("J134").Value = 1883
$ws.Range("K134").Value = 1871.64405
$ws.Range("L134").Value = 5649
$ws.Range("M134").Value = 663.3559500000001
$ws.Range("N134").Value = -10719
$ws.Range("H136").Value = 13987.632
$ws.Range("I136").Value = 838.7059
$ws.Range("K136").Value = 2516.1177
$ws.Range("M136").Value = 33.88229999999976

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2679.125
$ws.Range("I3").Value = 1800
$ws.Range("J3").Value = 8833
$ws.Range("K3").Value = 5400
$ws.Range("L3").Value = 26499
$ws.Range("M3").Value = -5288
$ws.Range("N3").Value = -26723
$ws.Range("H24").Value = 676.36365
$ws.Range("J24").Value = 780
$ws.Range("L24").Value = 2340
$ws.Range("N24").Value = -2800
$ws.Range("H59").Value = 2875
$ws.Range("I59").Value = 1000
$ws.Range("K59").Value = 3000
$ws.Range("M59").Value = -2460
$ws.Range("H74").Value = 10000
$ws.Range("J74").Value = 10000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -32122
$ws.Range("H77").Value = 10000
$ws.Range("J77").Value = 10000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -100608
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("N105").Value = 0
$ws.Range("H109").Value = 3828.5264
$ws.Range("I109").Value = 1328.5714
$ws.Range("J109").Value = 5286.8335
$ws.Range("K109").Value = 3985.7142
$ws.Range("L109").Value = 15860.5005
$ws.Range("M109").Value = -2945.7142
$ws.Range("N109").Value = -17940.5005
$ws.Range("H131").Value = 686.4299999999999
$ws.Range("J131").Value = 705.96704
$ws.Range("L131").Value = 2117.90112
$ws.Range("N131").Value = -12197.90112
$ws.Range("H138").Value = 1629.8823
$ws.Range("I138").Value = 1336.4286
$ws.Range("K138").Value = 4009.2858
$ws.Range("M138").Value = 1130.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 35718136
$ws.Range("I102").Value = 35718136
$ws.Range("K102").Value = 35718136
$ws.Range("M102").Value = -35716514

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3191.2632
$ws.Range("I40").Value = 2880.2856
$ws.Range("J40").Value = 4062
$ws.Range("K40").Value = 2880.2856
$ws.Range("L40").Value = 4062
$ws.Range("M40").Value = -2744.2856
$ws.Range("N40").Value = -4334
$ws.Range("H122").Value = 894135.3
$ws.Range("I122").Value = 1963387.8
$ws.Range("K122").Value = 5890163.4
$ws.Range("M122").Value = -5887713.4
$ws.Range("H136").Value = 879.3488
$ws.Range("I136").Value = 776.8461
$ws.Range("J136").Value = 1878.75
$ws.Range("K136").Value = 2330.5383
$ws.Range("L136").Value = 5636.25
$ws.Range("M136").Value = 219.4616999999998
$ws.Range("N136").Value = -10736.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 15000
$ws.Range("J54").Value = 15000
$ws.Range("L54").Value = 15000
$ws.Range("N54").Value = -16040
$ws.Range("H122").Value = 1336.6364
$ws.Range("I122").Value = 1352.4286
$ws.Range("J122").Value = 1005
$ws.Range("K122").Value = 4057.2858
$ws.Range("L122").Value = 3015
$ws.Range("M122").Value = -1607.2858
$ws.Range("N122").Value = -7915
